$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 362, shifting existing rows 362-453 down to 363-454
$ws.Rows.Item(362).Insert()

# Populate the newly inserted row 362 with the new weekly record
$ws.Cells.Item(362, 1).Value = 4
$ws.Cells.Item(362, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(362, 3).Value = "Los Lagos"
$ws.Cells.Item(362, 4).Value = 44722
$ws.Cells.Item(362, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(362, 5).Value = 10
$ws.Cells.Item(362, 6).Value = 100112006
$ws.Cells.Item(362, 7).Value = "Repollo"
$ws.Cells.Item(362, 8).Value = "Crespo record"
$ws.Cells.Item(362, 9).Value = "Primera"
$ws.Cells.Item(362, 10).Value = 1000
$ws.Cells.Item(362, 11).Value = 1700
$ws.Cells.Item(362, 12).Value = 1800
$ws.Cells.Item(362, 13).Value = 1750
$ws.Cells.Item(362, 14).Value = "`$/unidad"
$ws.Cells.Item(362, 15).Value = "Región Metropolitana"
$ws.Cells.Item(362, 16).Value = 1750
$ws.Cells.Item(362, 17).Value = 1
$ws.Cells.Item(362, 18).Value = "Hortaliza"
